$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily crypto price/volume refresh (GitHub Actions scrape, 2023-12-25).
# A leading "'" forces values that would otherwise look like plain numbers
# (e.g. "265.00") to stay text, matching the original "Price" column
# formatting instead of being auto-converted to a Number by Excel.

$ws.Range('D2').Value = '43.384.55'
$ws.Range('E2').Value = '  -0.88%  '

$ws.Range('D3').Value = '2.284.86'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''112.81'
$ws.Range('E5').Value = '  -1.28%  '

$ws.Range('D6').Value = '''265.00'
$ws.Range('E6').Value = '  -1.30%  '

$ws.Range('D7').Value = '''0.627'
$ws.Range('E7').Value = '  +0.55%  '

$ws.Range('E8').Value = '  +0.41%  '

$ws.Range('D9').Value = '''0.608'
$ws.Range('E9').Value = '  -1.92%  '

$ws.Range('D10').Value = '''47.08'
$ws.Range('E10').Value = '  -2.32%  '

$ws.Range('D11').Value = '''0.0935'
$ws.Range('E11').Value = '  -0.18%  '

$ws.Range('D12').Value = '''9.32'
$ws.Range('E12').Value = '  +6.41%  '

$ws.Range('E13').Value = '  +2.26%  '

$ws.Range('D14').Value = '''15.50'
$ws.Range('E14').Value = '  -0.68%  '

$ws.Range('D15').Value = '2.630.48'
$ws.Range('E15').Value = '  -0.14%  '

$ws.Range('D16').Value = '''0.865'
$ws.Range('E16').Value = '  +2.27%  '

$ws.Range('D17').Value = '2.282.72'
$ws.Range('E17').Value = '  -0.33%  '

$ws.Range('D18').Value = '43.284.87'
$ws.Range('E18').Value = '  -0.72%  '

$ws.Range('E19').Value = '  -1.08%  '

$ws.Range('D20').Value = '''6.80'
$ws.Range('E20').Value = '  +4.14%  '

$ws.Range('D21').Value = '''71.99'
$ws.Range('E21').Value = '  -0.66%  '

$ws.Range('E22').Value = '  -0.55%  '

$ws.Range('D23').Value = '''235.20'
$ws.Range('E23').Value = '  +0.95%  '

$ws.Range('D24').Value = '''9.51'
$ws.Range('E24').Value = '  -2.70%  '

$ws.Range('D25').Value = '''2.86'
$ws.Range('E25').Value = '  +1.70%  '

$ws.Range('E26').Value = '  +1.91%  '

$ws.Range('D27').Value = '''11.41'
$ws.Range('E27').Value = '  -1.34%  '

$ws.Range('D28').Value = '''41.02'
$ws.Range('E28').Value = '  -2.61%  '

$ws.Range('E29').Value = '  -1.66%  '

$ws.Range('E30').Value = '  -0.65%  '

$ws.Range('D31').Value = '''173.22'
$ws.Range('E31').Value = '  -1.92%  '

$ws.Range('D32').Value = '''21.56'
$ws.Range('E32').Value = '  -0.13%  '

$ws.Range('D33').Value = '''0.0906'
$ws.Range('E33').Value = '  -3.24%  '

$ws.Range('D34').Value = '''5.70'
$ws.Range('E34').Value = '  +2.64%  '

$ws.Range('E35').Value = '  +0.64%  '

$ws.Range('E36').Value = '  -1.90%  '

$ws.Range('D37').Value = '''0.0369'
$ws.Range('E37').Value = '  +3.66%  '

$ws.Range('D38').Value = '''4.04'
$ws.Range('E38').Value = '  +4.92%  '

$ws.Range('E39').Value = '  -4.20%  '

$ws.Range('D40').Value = '''2.62'
$ws.Range('E40').Value = '  +7.66%  '

$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').Value = '''76.22'
$ws.Range('E41').Value = '  +4.66%  '

$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').Value = '''14.05'
$ws.Range('E42').Value = '  +1.46%  '

$ws.Range('D43').Value = '''0.238'
$ws.Range('E43').Value = '  -2.02%  '

$ws.Range('D44').Value = '''6.15'
$ws.Range('E44').Value = '  +1.84%  '

$ws.Range('E45').Value = '  +0.24%  '

$ws.Range('D46').Value = '''1.38'
$ws.Range('E46').Value = '  -4.13%  '

$ws.Range('D47').Value = '''104.21'
$ws.Range('E47').Value = '  +1.55%  '

$ws.Range('E48').Value = '  -1.34%  '

$ws.Range('D49').Value = '''1.26'
$ws.Range('E49').Value = '  +2.97%  '

$ws.Range('D50').Value = '''0.0997'
$ws.Range('E50').Value = '  -0.54%  '

$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = '''0.600'
$ws.Range('E51').Value = '  +9.61%  '
